$wb = $excel.ActiveWorkbook

# Select the "Kappale" worksheet and add a new row of data (row 8)
$ws = $wb.Worksheets.Item("Kappale")

$ws.Range("A8").Value = "Numero"
$ws.Range("B8").Value = "kokonaisluku"
$ws.Range("D8").Value = "32bit"
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = "Kappaleen järjestys levyllä"

$ws.Range("E14").Select()
$ws.Activate()
